$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the comparison data for k=32 in columns E, F, G for rows 28-33
$ws.Range("E28").Value = 5
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 0.45

$ws.Range("E29").Value = 7
$ws.Range("F29").Value = 32
$ws.Range("G29").Value = 0.41

$ws.Range("E30").Value = 7
$ws.Range("F30").Value = 32
$ws.Range("G30").Value = 0.39

$ws.Range("E31").Value = 9
$ws.Range("F31").Value = 64
$ws.Range("G31").Value = 0.34

$ws.Range("E32").Value = 18
$ws.Range("F32").Value = 256
$ws.Range("G32").Value = 0.32

$ws.Range("E33").Value = 25
$ws.Range("F33").Value = 512
$ws.Range("G33").Value = 0.3

# Update the active cell selection to F34
$ws.Range("F34").Select()
